$d = $word.ActiveDocument

# Locate the anchor paragraph ("Charge that is separated is proportional to
# voltage") that the new content needs to follow, by scanning paragraphs for
# its text rather than assuming a fixed index.
$anchorIndex = -1
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "Charge that is separated is proportional to voltage*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Anchor paragraph not found"
}

$newLines = @(
  "B = (mu0*I/4*pi) * integral(dI*(r-r’)/modulus(r-r’)^3)",
  "Cyclotron - desktop device for accelerating particles",
  "Particles introduced into cyclotron",
  "Oscillating voltage within charges up particles",
  "Magnet within causes circulation of particles",
  "When current is alternated, particles accelerate faster and faster until eventually pops out through hole as very fast moving particle",
  "Can achieve 1000000 eV",
  "Every time particle reaches cyclotron frequency, particle gains boost in kinetic energy"
)

$insertAt = $anchorIndex
foreach ($line in $newLines) {
    $p = $d.Paragraphs.Item($insertAt)
    $p.Range.InsertParagraphAfter()
    $insertAt = $insertAt + 1
    $newPara = $d.Paragraphs.Item($insertAt)
    $newPara.Range.InsertAfter($line)
}

"inserted " + $newLines.Count + " paragraphs after index " + $anchorIndex
